# epexspot_prices.xlsx update
# 1) "Prix Spot" sheet: add column U ("04-jul") with header + 24 hourly values
# 2) "Gaz" sheet: add row 18 (2025-07-02 / 32.675)
# 3) "CO2" sheet: add row 18 (2025-07-02 / 71.40000000000001)

$wb = $excel.ActiveWorkbook

# --- 1) Prix Spot ---------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (T1) into U1 so the new
# header cell keeps the same bold/border/centered style, then overwrite
# its value with the new date label.
$wsSpot.Range("T1").Copy($wsSpot.Range("U1"))
$wsSpot.Range("U1").Value = "04-jul"

$spotValues = @(
    97.8,
    83,
    85.37,
    75.8,
    67.39,
    78.65000000000001,
    83.59,
    103.48,
    100,
    91.15000000000001,
    70.01000000000001,
    20.64,
    22.64,
    12.34,
    12.2,
    25.2,
    21.88,
    62.04,
    90.02,
    108.5,
    111.6,
    109.13,
    111.8,
    96.25
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 21).Value = $spotValues[$i]
}

# --- 2) Gaz ----------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date cell to be stored as text (matching the other "Date"
# column cells, which are plain inline/shared strings like "2025-06-16"),
# instead of letting Excel auto-convert the "yyyy-mm-dd" string into a
# date serial number. Resetting the style back to "Normal" afterwards
# clears the temporary text-number-format so the cell ends up unstyled,
# just like its neighbours.
$wsGaz.Range("A18").NumberFormat = "@"
$wsGaz.Range("A18").Value = "2025-07-02"
$wsGaz.Range("A18").Style = "Normal"
$wsGaz.Range("B18").Value = 32.675

# --- 3) CO2 ------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A18").NumberFormat = "@"
$wsCo2.Range("A18").Value = "2025-07-02"
$wsCo2.Range("A18").Style = "Normal"
$wsCo2.Range("B18").Value = 71.40000000000001
